$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig1AB")

# Thomas' feedback: rename the "BK/IAM" category label to "Model" in the
# Fig1AB data table (rows 2 and 3, column B).
$ws.Range("B2").Value = "Model"
$ws.Range("B3").Value = "Model"
